# Update Request Body structure in the data sheet (TestData)
# - Split the single "RequestBody" (JSON) column into two columns:
#   RequestBodyKey / RequestBodyValue
# - Shift ExpectedStatus / ValidationPath / ExpectedValue one column right
# - Re-point the mailto hyperlink that used to live in the old
#   ValidationPath-adjacent column to its new location
# - Adjust row heights / column widths / selection to match the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# --- 1. Remove existing hyperlinks up front (their ranges will move once we
#        insert a column, and this engine does not re-anchor them automatically) ---
$ws.Range("A1").Hyperlinks.Delete()

# --- 2. Insert a new column before the old "ExpectedStatus" column (G) so
#        the old RequestBody column (F) keeps its position and becomes
#        RequestBodyKey, while the brand new column G becomes
#        RequestBodyValue. Everything from the old G onward shifts right. ---
$ws.Columns("G").Insert()

# --- 3. Header row ---
$ws.Range("A1").Value = "TestID"
$ws.Range("B1").Value = "BaseURL"
$ws.Range("C1").Value = "Endpoint"
$ws.Range("D1").Value = "Method"
$ws.Range("E1").Value = "Headers"
$ws.Range("F1").Value = "RequestBodyKey"
$ws.Range("G1").Value = "RequestBodyValue"
$ws.Range("H1").Value = "ExpectedStatus"
$ws.Range("I1").Value = "ValidationPath"
$ws.Range("J1").Value = "ExpectedValue"

# --- 4. Data rows ---
# Row 2 - TC001
$ws.Range("A2").Value = "TC001"
$ws.Range("B2").Value = "https://jsonplaceholder.typicode.com"
$ws.Range("C2").Value = "/posts/1"
$ws.Range("D2").Value = "GET"
$ws.Range("E2").Value = "Content-Type: application/json"
$ws.Range("F2").Value = "null"
$ws.Range("G2").Value = "null"
$ws.Range("H2").Value = 200
$ws.Range("I2").Value = "$.title"
$ws.Range("J2").Value = "sunt aut facere"

# Row 3 - TC002
$ws.Range("A3").Value = "TC002"
$ws.Range("B3").Value = "https://reqres.in"
$ws.Range("C3").Value = "/api/users"
$ws.Range("D3").Value = "POST"
$ws.Range("E3").Value = "Content-Type: application/json"
$ws.Range("F3").Value = "name, job"
$ws.Range("G3").Value = "John Doe, QA Lead"
$ws.Range("H3").Value = 201
$ws.Range("I3").Value = "$.name"
$ws.Range("J3").Value = "John Doe"

# Row 4 - TC003
$ws.Range("A4").Value = "TC003"
$ws.Range("B4").Value = "https://reqres.in"
$ws.Range("C4").Value = "/api/users/2"
$ws.Range("D4").Value = "GET"
$ws.Range("E4").Value = "Content-Type: application/json"
$ws.Range("F4").Value = "name, job"
$ws.Range("G4").Value = "John Doe, QA Lead"
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = "$.data.email"
$ws.Range("J4").Value = "janet.weaver@reqres.in"

# Row 5 - TC004
$ws.Range("A5").Value = "TC004"
$ws.Range("B5").Value = "https://restcountries.com"
$ws.Range("C5").Value = "/v3.1/name/germany"
$ws.Range("D5").Value = "GET"
$ws.Range("E5").Value = "Content-Type: application/json"
$ws.Range("F5").Value = "tcid.id, dob.lastdob"
$ws.Range("G5").Value = "jj3, 77-09-9"
$ws.Range("H5").Value = 200
$ws.Range("I5").Value = "$.[0].capital[0]"
$ws.Range("J5").Value = "Berlin"

# Row 6 - TC005
$ws.Range("A6").Value = "TC005"
$ws.Range("B6").Value = "https://api.openweathermap.org"
$ws.Range("C6").Value = "/data/2.5/weather?q=London&appid=6aa6cd8c45d248d374aac371cd"
$ws.Range("D6").Value = "GET"
$ws.Range("E6").Value = "Content-Type: application/json"
$ws.Range("F6").Value = "null"
$ws.Range("G6").Value = "null"
$ws.Range("H6").Value = 200
$ws.Range("I6").Value = "$.name"
$ws.Range("J6").Value = "London"

# Row 7 - TC006
$ws.Range("A7").Value = "TC006"
$ws.Range("B7").Value = "https://api.github.com"
$ws.Range("C7").Value = "/repos/octocat/Hello-World"
$ws.Range("D7").Value = "GET"
$ws.Range("E7").Value = "Accept: application/vnd.github.v3+json"
$ws.Range("F7").Value = "null"
$ws.Range("G7").Value = "null"
$ws.Range("H7").Value = 200
$ws.Range("I7").Value = "$.name"
$ws.Range("J7").Value = "Hello-World"

# --- 5. New blank, styled row far below the table (matches the author's
#        stray formatted cells at C15:D15) ---
$ws.Range("C2").Copy()
$ws.Range("C15:D15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 6. Re-create the hyperlinks at their (possibly shifted) locations ---
# NOTE: Hyperlinks.Add overwrites the cell's text with TextToDisplay, so that
# argument must be the original cell value (not the link's display/address).
$ws.Hyperlinks.Add($ws.Range("B2"), "https://jsonplaceholder.typicode.com/", "", "", "https://jsonplaceholder.typicode.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://reqres.in/", "", "", "https://reqres.in")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://reqres.in/", "", "", "https://reqres.in")
$ws.Hyperlinks.Add($ws.Range("J4"), "mailto:janet.weaver@reqres.in", "", "", "janet.weaver@reqres.in")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://restcountries.com/", "", "", "https://restcountries.com")
$ws.Hyperlinks.Add($ws.Range("B6"), "https://api.openweathermap.org/", "", "", "https://api.openweathermap.org")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://api.github.com/", "", "", "https://api.github.com")

# --- 7. Row heights to match the new (narrower) columns ---
$ws.Rows(1).RowHeight = 30
$ws.Rows(2).RowHeight = 30
$ws.Rows(3).RowHeight = 30
$ws.Rows(4).RowHeight = 45
$ws.Rows(5).RowHeight = 30
$ws.Rows(6).RowHeight = 30
$ws.Rows(7).RowHeight = 45

# --- 8. Column widths: E/F/G get explicit widths, H-J fall back to default ---
$ws.Range("E1").EntireColumn.ColumnWidth = 19.59
$ws.Range("F1").EntireColumn.ColumnWidth = 26.74
$ws.Range("G1").EntireColumn.ColumnWidth = 26.74
$ws.Range("J1").EntireColumn.ColumnWidth = 8.43

# --- 9. Selection matches the author's final cursor position ---
$ws.Range("G6").Select()
